# Applies the "丽水-漫展信息" gh-pages data refresh (commit 456a3b4) to both
# the "展览" and "全部类型" sheets, which carry identical event tables.
#
# Summary of the change:
#   - Row 2 (C2..I2) is rewritten from the "thp01" event to the "HP国风"
#     carnival event; B2 (date "2024-07-27") is untouched.
#   - Row 3 (B3..I3) is rewritten from the "HP国风" carnival event back to
#     the "thp01" event, one day later.
#   - Row 6: only the end time (E6) and the "want to go" count (F6) change.
#   - A brand-new row is inserted at position 7 ("丽水·R动漫嘉年华"), pushing
#     the former row 7 (LZ栗子) down to row 8.

# Writes a plain text value into a cell without letting "YYYY-MM-DD"-looking
# strings get auto-recognised (and reformatted) as real dates: force Text
# number format for the assignment, then restore the cell's normal
# (style-less) look by pasting formats from a known plain donor cell, so
# the end result matches the rest of the column (plain inline/shared string,
# no explicit style).
function Set-TextValue {
    param($ws, [string]$addr, [string]$value, $donor)

    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $value
    $donor.Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Application.CutCopyMode = 0
}

function Set-EventRow {
    param(
        $ws,
        [int]$r,
        $donor,
        [string]$b,
        [string]$c,
        [string]$d,
        [string]$e,
        $f,
        $g,
        [string]$h,
        [string]$i
    )
    if ($b -ne $null) {
        Set-TextValue $ws "B$r" $b $donor
    }
    $ws.Range("C$r").Value = $c
    $ws.Range("D$r").Value = $d
    $ws.Range("E$r").Value = $e
    $ws.Range("F$r").Value = $f
    $ws.Range("G$r").Value = $g
    $ws.Range("H$r").Value = $h
    $ws.Range("I$r").Value = $i
}

function Update-SheetData {
    param($ws)

    # A plain, never-styled cell (sheet header "Link") used purely as a
    # formatting donor for the Set-TextValue trick above.
    $donor = $ws.Range("H1")

    # --- Push the old row 7 ("丽水·LZ栗子动漫游戏嘉年华") down to row 8,
    #     values + formatting together, before row 7 gets overwritten with
    #     the new "丽水·R动漫嘉年华" event below. Using Copy(Destination)
    #     (rather than Rows.Insert, which invents a fresh/partial auto
    #     style for the vacated row) keeps column A's existing "index"
    #     style intact with no stray style entries.
    $ws.Range("A7:I7").Copy($ws.Range("A8:I8"))

    # --- Row 2: was "丽水·thp01～风摄少微", now the HP carnival entry.
    #     B2 ("2024-07-27") is not part of this change.
    Set-EventRow $ws 2 $donor $null "丽水·第四届HP国风动漫游戏嘉年华" "城北街798号 莱茵体育生活馆" "2024.07.27 08:30-07.27 17:00" 501 65 "https://show.bilibili.com/platform/detail.html?id=87305" "//i2.hdslb.com/bfs/openplatform/202406/YUnPOKGV1718268952725.jpeg"

    # --- Row 3: was the HP carnival entry, now "丽水·thp01～风摄少微" (date +1 day).
    Set-EventRow $ws 3 $donor "2024-07-28" "丽水·thp01～风摄少微" "大猷街 应星楼" "2024.07.28 10:00-07.28 18:00" 31 50 "https://show.bilibili.com/platform/detail.html?id=87134" "//i2.hdslb.com/bfs/openplatform/202407/WbSdFFLd1721636456044.jpeg"

    # --- Row 6: AEO event end time + want-to-go count update only.
    $ws.Range("E6").Value = "2024.08.17 09:00-08.17 18:00"
    $ws.Range("F6").Value = 690

    # --- Row 7 (new): 丽水·R动漫嘉年华 (A7's index value 6 is left as-is,
    #     matching the diff, which leaves <c r="A7"><v>6</v></c> untouched).
    Set-EventRow $ws 7 $donor "2024-08-24" "丽水·R动漫嘉年华" "中东路848号(解放街交汇) 飞达国际大酒店" "2024.08.24 09:30-08.25 17:00" 0 45 "https://show.bilibili.com/platform/detail.html?id=89651" "//i0.hdslb.com/bfs/openplatform/202407/7o5ALbAM1721383424201.jpeg"

    # Row 8 already carries the correct (former row-7) data after the copy
    # above; only its index cell needs bumping from 6 to 7.
    $ws.Range("A8").Value = 7
}

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
Update-SheetData $ws1

$ws4 = $wb.Worksheets.Item("全部类型")
Update-SheetData $ws4
